# Update Reaction_number (column C) values on the "NBR" and "BAR" sheets
# per the commit "add new ecoli studies and NB ratio analysis".

$wb = $excel.ActiveWorkbook

$nbrValues = @(697, 683, 679, 672, 660, 656, 650, 649, 649, 642, 638, 628, 622, 621, 616, 616, 614, 611, 606)
$barValues = @(651, 649, 650, 649, 650, 651, 650, 658, 656, 656, 654, 665, 662, 651, 650, 651, 651, 650, 650)

$wsNbr = $wb.Worksheets.Item("NBR")
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNbr.Cells.Item($row, 3).Value = $nbrValues[$i]
}

$wsBar = $wb.Worksheets.Item("BAR")
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBar.Cells.Item($row, 3).Value = $barValues[$i]
}
